$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the "EK" column header in C2 to "Buy-Price"
$ws.Range("C2").Value = "Buy-Price"

# Restore the selection that was active when the file was saved
$ws.Range("E19").Select()
